# Leitung als Spalte in Versuchspersonen nach Einrichtung ergaenzt
#
# Inserts a new "Leitung" column between the existing "Einrichtung" column (A)
# and the "SBBZ" column (B), shifting SBBZ/HZE one column to the right
# (B->C, C->D), and fills in the new column's header + per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns B (SBBZ) and C (HZE) one to the right by inserting
# a new, empty column at B.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Leitung"

# Values for the new "Leitung" column, row by row (rows 2-23).
$leitungValues = @(3, 6, 6, 12, 4, 2, 7, 6, 2, 1, 2, 5, 0, 2, 6, 3, 3, 2, 8, 5, 1, 2)

$row = 2
foreach ($value in $leitungValues) {
    $ws.Cells.Item($row, 2).Value = $value
    $row = $row + 1
}
